$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.993.90"
$ws.Range("E2").Value = "  -6.28%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.073.94"
$ws.Range("E3").Value = "  -9.04%  "

$ws.Range("E4").Value = "  +0.20%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "502.09"
$ws.Range("E5").Value = "  -5.01%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "163.65"
$ws.Range("E6").Value = "  -12.99%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.572"
$ws.Range("E7").Value = "  -5.21%  "

$ws.Range("E8").Value = "  +0.12%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.068.91"
$ws.Range("E9").Value = "  -9.12%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.571"
$ws.Range("E10").Value = "  -8.90%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "50.56"
$ws.Range("E11").Value = "  -13.95%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.124"
$ws.Range("E12").Value = "  -7.33%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000239"
$ws.Range("E13").Value = "  -6.05%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.51"
$ws.Range("E14").Value = "  -8.13%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.601.48"
$ws.Range("E15").Value = "  -8.06%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.111"
$ws.Range("E16").Value = "  -9.62%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.103.08"
$ws.Range("E17").Value = "  -8.13%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "61.049.26"
$ws.Range("E18").Value = "  -5.99%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.38"
$ws.Range("E19").Value = "  -6.08%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.34"
$ws.Range("E20").Value = "  -6.98%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.914"
$ws.Range("E21").Value = "  -6.00%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "349.10"
$ws.Range("E22").Value = "  -6.06%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "77.43"
$ws.Range("E23").Value = "  -4.90%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.55"
$ws.Range("E24").Value = "  -4.57%  "

$ws.Range("B25").Value = "LEO"
$ws.Range("C25").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "6.06"
$ws.Range("E25").Value = "  +4.31%  "

$ws.Range("B26").Value = "RenderToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.46"
$ws.Range("E26").Value = "  -3.49%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.73"
$ws.Range("E27").Value = "  -0.64%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.49"
$ws.Range("E28").Value = "  -6.01%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.59"
$ws.Range("E29").Value = "  -8.03%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.71"
$ws.Range("E30").Value = "  -9.64%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "622.89"
$ws.Range("E31").Value = "  -7.01%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "26.88"
$ws.Range("E32").Value = "  -9.13%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.04"
$ws.Range("E33").Value = "  -10.30%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "10.73"
$ws.Range("E34").Value = "  -4.26%  "

$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "57.93"
$ws.Range("E35").Value = "  -5.25%  "

$ws.Range("B36").Value = "Dai"
$ws.Range("C36").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.01%  "

$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0985"
$ws.Range("E37").Value = "  -6.81%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "34.57"
$ws.Range("E38").Value = "  -5.56%  "

$ws.Range("B39").Value = "FirstDigitalUSD"
$ws.Range("C39").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  +0.43%  "

$ws.Range("B40").Value = "TheGraph"
$ws.Range("C40").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.356"
$ws.Range("E40").Value = "  -6.24%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0₃0646"
$ws.Range("E41").Value = "  +4.00%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.117"
$ws.Range("E42").Value = "  -7.85%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.743.84"
$ws.Range("E43").Value = "  -3.85%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.36"
$ws.Range("E44").Value = "  +0.85%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.57"
$ws.Range("E45").Value = "  -3.20%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.74"
$ws.Range("E46").Value = "  +6.67%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0368"
$ws.Range("E47").Value = "  -6.25%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.42"
$ws.Range("E48").Value = "  -12.01%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.83"
$ws.Range("E49").Value = "  -0.24%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.118"
$ws.Range("E50").Value = "  -5.41%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "128.36"
$ws.Range("E51").Value = "  -7.00%  "
